$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Step 1: remove the <w:proofErr .../> wrapper around the "Sentinel" run
# (paragraph 4: the title line "Sentinel")
# ------------------------------------------------------------------
$pSentinel = $d.Paragraphs.Item(4)
$xmlSentinel = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:sz w:val="72"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="72"/></w:rPr><w:t>Sentinel</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pSentinel.Range.InsertXML($xmlSentinel)

# ------------------------------------------------------------------
# Step 2: remove the <w:lastRenderedPageBreak/> before "Installation"
# (paragraph 6: the "Installation" heading)
# ------------------------------------------------------------------
$pInstallation = $d.Paragraphs.Item(6)
$xmlInstallation = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="berschrift1"/></w:pPr><w:r><w:t>Installation</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pInstallation.Range.InsertXML($xmlInstallation)

# ------------------------------------------------------------------
# Step 3: add a new "Konfiguration" list item after "Aktualisierung"
# (paragraph 8) and before the following page-break paragraph
# (paragraph 9), continuing the same numbered/bulleted list.
# ------------------------------------------------------------------
$pPageBreakBeforeBenutzung = $d.Paragraphs.Item(9)
$pPageBreakBeforeBenutzung.Range.InsertParagraphBefore()
$pKonfiguration = $d.Paragraphs.Item(9)
$pKonfiguration.Range.Text = "Konfiguration"
$pKonfiguration.Style = "Listenabsatz"
$listTemplate = $d.Paragraphs.Item(8).Range.ListFormat.ListTemplate
$pKonfiguration.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true, 2, $false, $false)

# ------------------------------------------------------------------
# Step 4: append a page break and the new "Server konfiguration"
# paragraphs after "Benutzungshandbuch" (now paragraph 11), and
# drop the lastRenderedPageBreak/_GoBack bookmark on that heading.
# ------------------------------------------------------------------
$pBenutzung = $d.Paragraphs.Item(11)

# Reserve 7 empty paragraphs (1 page break + 6 text lines) up front so
# none of the paragraphs we fill in are ever the document's very last
# paragraph while we set their content (avoids a spurious trailing
# empty paragraph being introduced).
for ($k = 0; $k -lt 7; $k++) {
    $d.Paragraphs.Item($d.Paragraphs.Count).Range.InsertParagraphAfter()
}

$pPageBreak2 = $d.Paragraphs.Item(12)
$xmlPageBreak2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:br w:type="page"/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pPageBreak2.Range.InsertXML($xmlPageBreak2)

$newTexts = @(
    "Server konfiguration:",
    "Dem Server kann man mehrere Programm Argumente mitgeben:",
    "-debugmode",
    "-port",
    "-ip",
    "-headless"
)
for ($k = 0; $k -lt $newTexts.Length; $k++) {
    $idx = 13 + $k
    $para = $d.Paragraphs.Item($idx)
    $para.Style = "Standard"
    $para.Range.Text = $newTexts[$k]
}

# Now that "Benutzungshandbuch" is no longer the last paragraph in the
# document, clean its lastRenderedPageBreak and the _GoBack bookmark.
$xmlBenutzung = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="berschrift1"/></w:pPr><w:r><w:t>Benutzungshandbuch</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pBenutzung.Range.InsertXML($xmlBenutzung)
